# fix confusion between integer and number in schema_from_excel
#
# The bug: a "number" column whose values all happened to look like
# integers (32, 25, 12 ...) was being inferred as an "integer" column by
# schema_from_excel. table1!C3 is changed from the integer 32 to the
# non-integer 32.1 so the column is unambiguously a "number" in the test
# fixture. The rest of the diff is just editor/window state left over
# from opening the file to make the edit (active sheet/cell, and the
# page setup that got touched on table3 along the way).

$wb = $excel.ActiveWorkbook

# --- table1 (sheet1): the actual data fix -------------------------------
$ws1 = $wb.Worksheets.Item("table1")
$ws1.Range("C3").Value = 32.1

# --- table2 (sheet2): selection moved, no longer the active tab --------
$ws2 = $wb.Worksheets.Item("table2")
[void]$ws2.Activate()
[void]$ws2.Range("K10").Select()

# --- table3 (sheet3): selection moved, and picked up a page setup ------
$ws3 = $wb.Worksheets.Item("table3")
[void]$ws3.Activate()
[void]$ws3.Range("D15").Select()
$ws3.PageSetup.Orientation = 1

# --- table1 (sheet1) ends up the active sheet/tab, selection on C4 -----
[void]$ws1.Activate()
[void]$ws1.Range("C4").Select()
